$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.gender = `"FEMALE`"`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$newStatQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.gender = `"FEMALE`"`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

$ws.Range("B2").Value = $newQuery
$ws.Range("C2").Value = $newStatQuery

$ws.Rows.Item(2).RowHeight = 174

$ws.Range("B6").Select()
